# Cross check all the test script.
# Adds a new "UnverifiedVendor" worksheet (with unverified-vendor data) between
# the "user" and "Customer" sheets, and makes "InvalidData" the active tab.
#
# NOTE: worksheet object references returned by this host resolve by
# position, not by a stable identity -- once any sheet is inserted/deleted,
# previously-captured worksheet variables can start pointing at a different
# sheet. To stay safe we always re-fetch a worksheet **by name** immediately
# before using it, right after any Add()/Delete() call.

$wb = $excel.ActiveWorkbook

# --- Bump the internal sheetId counter up to 5 -----------------------------
# Excel assigns each new sheet a sheetId of (current max sheetId + 1). The
# existing sheets use sheetId 1 (user), 3 (Customer), 2 (InvalidData) -- max
# is 3, so a single Add() would land on sheetId 4. The target file expects
# the new sheet to carry sheetId 5, so we add a throw-away sheet first (it
# claims sheetId 4), add our real sheet right after "user" (claims sheetId
# 5), then delete the throw-away sheet.
$scratch = $wb.Worksheets.Add()
$userSheet = $wb.Worksheets.Item("user")
$wb.Worksheets.Add([System.Reflection.Missing]::Value, $userSheet) | Out-Null

$scratch = $wb.Worksheets.Item("Sheet1")
$scratch.Delete()

$newSheet = $wb.Worksheets.Item("Sheet2")
$newSheet.Name = "UnverifiedVendor"

# --- Header row ------------------------------------------------------------
$newSheet.Range("A1").Value = "MobileNumber"
$newSheet.Range("B1").Value = "Type"
$newSheet.Range("C1").Value = "EmailAddress"
$newSheet.Range("D1").Value = "Password"
$newSheet.Range("E1").Value = "Location"
$newSheet.Range("F1").Value = "IsEmailVerified"
$newSheet.Range("G1").Value = "IsMobileVerified"
$newSheet.Range("H1").Value = "IsAadharVerified"

# --- Row 2 -------------------------------------------------------------
$newSheet.Range("A2").Value = 7785683689
$newSheet.Range("B2").Value = "v"
$newSheet.Range("C2").Value = "j0tx13058c@mxscout.com"
$newSheet.Range("D2").Value = "Aquaclean@123"
$newSheet.Range("E2").Value = "Pune"
$newSheet.Range("F2").Value = 1
$newSheet.Range("G2").Value = 1
$newSheet.Range("H2").Value = 0

# --- Row 3 -------------------------------------------------------------
$newSheet.Range("A3").Value = 9846789648
$newSheet.Range("B3").Value = "v"
$newSheet.Range("C3").Value = "isha342@yahoo.com"
$newSheet.Range("D3").Value = "Isha@123"
$newSheet.Range("E3").Value = "Pune"
$newSheet.Range("F3").Value = 1
$newSheet.Range("G3").Value = 0
$newSheet.Range("H3").Value = 1

# --- Row 4 -------------------------------------------------------------
$newSheet.Range("A4").Value = 9645132789
$newSheet.Range("B4").Value = "v"
$newSheet.Range("C4").Value = "ishwar94@gmail.com"
$newSheet.Range("D4").Value = "Ishwar@123"
$newSheet.Range("E4").Value = "Pune"
$newSheet.Range("F4").Value = 0
$newSheet.Range("G4").Value = 1
$newSheet.Range("H4").Value = 1

# --- Cell formatting ---------------------------------------------------
# A2 & C2: left + top aligned
$newSheet.Range("A2").HorizontalAlignment = -4131
$newSheet.Range("A2").VerticalAlignment = -4160
$newSheet.Range("C2").HorizontalAlignment = -4131
$newSheet.Range("C2").VerticalAlignment = -4160

# A3 / A4: left aligned
$newSheet.Range("A3").HorizontalAlignment = -4131
$newSheet.Range("A4").HorizontalAlignment = -4131

# D2:D4 stored as text
$newSheet.Range("D2:D4").NumberFormat = "@"

# --- Hyperlink on C2 -----------------------------------------------------
$newSheet.Hyperlinks.Add($newSheet.Range("C2"), "mailto:j0tx13058c@mxscout.com") | Out-Null

# --- Column widths (match source best-fit layout) -------------------------
$newSheet.Columns.Item(1).ColumnWidth = 13.5546875
$newSheet.Columns.Item(2).ColumnWidth = 5.33203125
$newSheet.Columns.Item(3).ColumnWidth = 24.109375
$newSheet.Columns.Item(4).ColumnWidth = 15.5546875
$newSheet.Columns.Item(5).ColumnWidth = 8.44140625
$newSheet.Columns.Item(6).ColumnWidth = 14.6640625
$newSheet.Columns.Item(7).ColumnWidth = 16.109375
$newSheet.Columns.Item(8).ColumnWidth = 16.109375

# --- Selection on the new sheet -------------------------------------------
$newSheet.Range("C29").Select()

# --- Selection on "user" sheet changes too --------------------------------
$userSheet = $wb.Worksheets.Item("user")
$userSheet.Range("A6").Select()

# --- Make InvalidData the active / selected tab ---------------------------
$invalid = $wb.Worksheets.Item("InvalidData")
$invalid.Activate()
$invalid.Range("E2").Select()
